$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "target" column (G) values were relabelled from "deuteron" to "d"
# for every data row.
$ws.Range("G2:G13").Value = "d"

# The header row (A1:K1) was given a bold, centered style.
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# Reflect the new selection over the (now styled) header row.
[void]$header.Select()
